# "Added saved and applied filter test"
#
# - Renames TC041 -> TC08
# - Adds a RunAsLabel column to TC02
# - Adds four new test-case sheets: TC09, TC10, TC11, TC12
# - Tweaks a couple of selections
# - Makes TC11 the active tab

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------
# 1. Rename TC041 -> TC08 (sheetId / r:id stay put on rename)
# ---------------------------------------------------------------
$tc08 = $wb.Worksheets.Item("TC041")
$tc08.Name = "TC08"

# ---------------------------------------------------------------
# 2. Insert a "RunAsLabel" column into TC02 (new column E)
# ---------------------------------------------------------------
$tc02 = $wb.Worksheets.Item("TC02")
$tc02.Activate()
$tc02.Columns("E:E").Insert()
$tc02.Range("E1").Value = "RunAsLabel"
$tc02.Range("E2").Value = "aaron.rucker@uthsc.edu"
$tc02.Range("E3").Value = "CWC\Admin"
$tc02.Range("E4").Value = "JORDAN\User"
$tc02.Range("E1").Interior.Color = 65535
$tc02.Columns("E:E").ColumnWidth = 25.85546875
$tc02.Range("E4").Select()

# ---------------------------------------------------------------
# 3. Tidy up TC03's selection (no data change)
# ---------------------------------------------------------------
$tc03 = $wb.Worksheets.Item("TC03")
$tc03.Activate()
$tc03.Range("A1:D2").Select()

# ---------------------------------------------------------------
# 4. Add four new sheets after TC051.
#    A throwaway sheet is created (and removed) first purely to
#    burn a sheetId, so TC09 lands on sheetId 14 (matching the
#    target workbook.xml) instead of 13.
# ---------------------------------------------------------------
$tc051 = $wb.Worksheets.Item("TC051")
$placeholder = $wb.Worksheets.Add($null, $tc051)

$tc09 = $wb.Worksheets.Add($null, $placeholder)
$tc09.Name = "TC09"
$placeholder.Delete()

$tc10 = $wb.Worksheets.Add($null, $tc09)
$tc10.Name = "TC10"

$tc11 = $wb.Worksheets.Add($null, $tc10)
$tc11.Name = "TC11"

$tc12 = $wb.Worksheets.Add($null, $tc11)
$tc12.Name = "TC12"

# ---------------------------------------------------------------
# 5. TC09 - ChangeStatus test data
#    Shared-string order: E1, E2, D1, D2, then the rest.
# ---------------------------------------------------------------
$tc09.Range("E1").Value = "ChangeStatus"
$tc09.Range("E2").Value = "password changed successfully"
$tc09.Range("D1").Value = "NewPassword"
$tc09.Range("D2").Value = "Test1"
$tc09.Range("A1").Value = "username"
$tc09.Range("B1").Value = "password"
$tc09.Range("C1").Value = "waitfor"
$tc09.Range("A2").Value = "NextGen\Nagendra"
$tc09.Range("B2").Value = "Test"
$tc09.Range("C2").Value = 10

$tc09.Range("A1:E1").Interior.Color = 65535
$tc09.Range("B1").BorderAround(1)
$tc09.Range("A1:E2").Select()

# ---------------------------------------------------------------
# 6. TC10 - NotAnumber/NaN test data
#    Shared-string order: D1, D2.
# ---------------------------------------------------------------
$tc10.Range("D1").Value = "NotAnumber"
$tc10.Range("D2").Value = "NaN"
$tc10.Range("A1").Value = "username"
$tc10.Range("B1").Value = "password"
$tc10.Range("C1").Value = "RunAs"
$tc10.Range("E1").Value = "waitFor"
$tc10.Range("A2").Value = "NextGen\Nagendra"
$tc10.Range("B2").Value = "WSCAdmin"
$tc10.Range("C2").Value = "UTMG"
$tc10.Range("E2").Value = 10

$tc10.Range("A1:E1").Interior.Color = 65535
$tc10.Range("B1").BorderAround(1)
$tc10.Columns("A:A").ColumnWidth = 24.28515625
$tc10.Columns("B:D").ColumnWidth = 30.7109375
$tc10.Range("A1:E2").Select()

# ---------------------------------------------------------------
# 7. TC11 - Applied/Saved filter test data (the headline feature)
#    Shared-string order: D1, D2, E1, F1, E2, F2.
# ---------------------------------------------------------------
$tc11.Range("D1").Value = "Practices"
$tc11.Range("D2").Value = "University Of Tennessee,UT Medical Group, Inc."
$tc11.Range("E1").Value = "AppliedFilterTitle"
$tc11.Range("F1").Value = "SavedFilterTitle"
$tc11.Range("E2").Value = "Applied Filter"
$tc11.Range("F2").Value = "Saved Filter"
$tc11.Range("A1").Value = "username"
$tc11.Range("B1").Value = "password"
$tc11.Range("C1").Value = "RunAs"
$tc11.Range("G1").Value = "waitFor"
$tc11.Range("A2").Value = "NextGen\Nagendra"
$tc11.Range("B2").Value = "WSCAdmin"
$tc11.Range("C2").Value = "UTMG"
$tc11.Range("G2").Value = 10

$tc11.Range("A1:G1").Interior.Color = 65535
$tc11.Range("B1").BorderAround(1)
$tc11.Columns("A:C").ColumnWidth = 30.7109375
$tc11.Columns("D:F").ColumnWidth = 42.85546875
$tc11.Columns("G:G").ColumnWidth = 30.7109375
$tc11.Range("A2").Select()

# ---------------------------------------------------------------
# 8. TC12 - plain RunAs smoke test, no new strings
# ---------------------------------------------------------------
$tc12.Range("A1").Value = "username"
$tc12.Range("B1").Value = "password"
$tc12.Range("C1").Value = "RunAs"
$tc12.Range("D1").Value = "waitFor"
$tc12.Range("A2").Value = "NextGen\Nagendra"
$tc12.Range("B2").Value = "WSCAdmin"
$tc12.Range("C2").Value = "UTMG"
$tc12.Range("D2").Value = 10

$tc12.Range("A1:D1").Interior.Color = 65535
$tc12.Range("B1").BorderAround(1)
$tc12.Columns("A:D").ColumnWidth = 30.7109375
$tc12.Range("D1").Select()

# ---------------------------------------------------------------
# 9. TC11 becomes the active tab (and bookViews activeTab moves
#    with it automatically).
# ---------------------------------------------------------------
$tc11.Activate()
